$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('Sergipe','Estuda e trabalha','31/12/2012',11.96582770060558),
    @('Sergipe','Estuda e trabalha','31/12/2013',11.97356237363959),
    @('Sergipe','Estuda e trabalha','31/12/2014',11.70599199124013),
    @('Sergipe','Estuda e trabalha','31/12/2015',10.84071243896443),
    @('Sergipe','Estuda e trabalha','31/12/2016',9.03226368509583),
    @('Sergipe','Estuda e trabalha','31/12/2017',9.652746847417564),
    @('Sergipe','Estuda e trabalha','31/12/2018',8.001179706823077),
    @('Sergipe','Estuda e trabalha','31/12/2019',9.028736239756153),
    @('Sergipe','Estuda e trabalha','31/12/2020',6.599897911123117),
    @('Sergipe','Estuda e trabalha','31/12/2021',7.787070507450794),
    @('Sergipe','Estuda e trabalha','31/12/2022',8.859857473132625),
    @('Sergipe','Estuda e trabalha','31/12/2023',9.143685441195133),
    @('Sergipe','Estuda e trabalha','31/12/2024',9.682364068073488),
    @('Sergipe','Não estuda e não trabalha','31/12/2012',24.22340546063874),
    @('Sergipe','Não estuda e não trabalha','31/12/2013',24.52260832718279),
    @('Sergipe','Não estuda e não trabalha','31/12/2014',23.63404389166559),
    @('Sergipe','Não estuda e não trabalha','31/12/2015',25.55528314284463),
    @('Sergipe','Não estuda e não trabalha','31/12/2016',27.87741724975137),
    @('Sergipe','Não estuda e não trabalha','31/12/2017',30.73187980949513),
    @('Sergipe','Não estuda e não trabalha','31/12/2018',31.34978958607622),
    @('Sergipe','Não estuda e não trabalha','31/12/2019',27.53425958847413),
    @('Sergipe','Não estuda e não trabalha','31/12/2020',32.31050202002471),
    @('Sergipe','Não estuda e não trabalha','31/12/2021',32.97505379805603),
    @('Sergipe','Não estuda e não trabalha','31/12/2022',27.11898162868759),
    @('Sergipe','Não estuda e não trabalha','31/12/2023',26.87570132222146),
    @('Sergipe','Não estuda e não trabalha','31/12/2024',24.32929539320264),
    @('Sergipe','Só estuda','31/12/2012',25.62430701642132),
    @('Sergipe','Só estuda','31/12/2013',26.68781488040046),
    @('Sergipe','Só estuda','31/12/2014',24.98315340218343),
    @('Sergipe','Só estuda','31/12/2015',25.81813976740516),
    @('Sergipe','Só estuda','31/12/2016',25.93601653341085),
    @('Sergipe','Só estuda','31/12/2017',25.0273039182603),
    @('Sergipe','Só estuda','31/12/2018',28.87947423748897),
    @('Sergipe','Só estuda','31/12/2019',26.48073497960748),
    @('Sergipe','Só estuda','31/12/2020',30.0133146393255),
    @('Sergipe','Só estuda','31/12/2021',26.88217833376401),
    @('Sergipe','Só estuda','31/12/2022',29.07402643570408),
    @('Sergipe','Só estuda','31/12/2023',27.08038340032537),
    @('Sergipe','Só estuda','31/12/2024',24.2705199184941),
    @('Sergipe','Só trabalha','31/12/2012',38.18645982233473),
    @('Sergipe','Só trabalha','31/12/2013',36.81601441877716),
    @('Sergipe','Só trabalha','31/12/2014',39.6768107149108),
    @('Sergipe','Só trabalha','31/12/2015',37.78586465078588),
    @('Sergipe','Só trabalha','31/12/2016',37.15430253174191),
    @('Sergipe','Só trabalha','31/12/2017',34.58806942482718),
    @('Sergipe','Só trabalha','31/12/2018',31.7695564696118),
    @('Sergipe','Só trabalha','31/12/2019',36.95626919216218),
    @('Sergipe','Só trabalha','31/12/2020',31.07628542952682),
    @('Sergipe','Só trabalha','31/12/2021',32.35569736072931),
    @('Sergipe','Só trabalha','31/12/2022',34.94713446247558),
    @('Sergipe','Só trabalha','31/12/2023',36.90022983625806),
    @('Sergipe','Só trabalha','31/12/2024',41.71782062022998)
)

$numRows = $data.Count
$arr = New-Object 'object[,]' $numRows,4
for ($i = 0; $i -lt $numRows; $i++) {
    $arr[$i,0] = $data[$i][0]
    $arr[$i,1] = $data[$i][1]
    $arr[$i,2] = $data[$i][2]
    $arr[$i,3] = $data[$i][3]
}

$startRow = 2
$endRow = $startRow + $numRows - 1
$rng = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 4))
$rng.Value = $arr
